$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = -3
$ws.Range("F6").Value = 5
$ws.Range("F10").Value = -4
$ws.Range("F11").Value = -8
$ws.Range("F12").Value = -3
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = -2
$ws.Range("F19").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("F25").Value = 6
$ws.Range("F26").Value = 0
